# Automatische test-sync: 2025-08-05 18:08:50
#
# Adds the new "Testmail #2" log entry (row 23) to the "Logs" sheet,
# extends the conditional-formatting ranges that covered rows 2-22 so
# they also cover the new row 2-23, and bumps the "Planning / Afspraak"
# tally on the "Dashboard" sheet from 16 to 17.

$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append the new row ------------------------------------
$ws = $wb.Worksheets.Item("Logs")

$newRow = 23
$ws.Cells.Item($newRow, 1).Value  = "Wil je dit oppakken?"
$ws.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value  = "Testmail #2: Wil je dit oppakken?"
$ws.Cells.Item($newRow, 4).Value  = "Planning / Afspraak"
$ws.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Cells.Item($newRow, 6).Value  = "2025-08-05 18:08:46"
$ws.Cells.Item($newRow, 7).Value  = "Ja"
$ws.Cells.Item($newRow, 8).Value  = "Ja"
$ws.Cells.Item($newRow, 9).Value  = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# --- extend conditional formatting sqref D/G/H/I/J 2:22 -> 2:23 ----------
$ws.Range("D2:D22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D23"))
$ws.Range("G2:G22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G23"))
$ws.Range("H2:H22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H23"))
$ws.Range("I2:I22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I23"))
$ws.Range("J2:J22").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J23"))

# --- "Dashboard" sheet: bump the Planning / Afspraak count ---------------
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("B2").Value = 17
